$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in the header (F1)
$ws.Range("F1").Value = "Last status check on: 16.02.2022 06:30"

# Makro row (row 5): Delta Cena (D5) becomes a real number instead of text "+0.4"
$ws.Range("D5").Value = 0.4

# Makro row (row 5): Old Datum (E5) becomes a real Excel date-time serial
# instead of plain text, formatted the same way as the other rows in column E
$ws.Range("E5").Value = 44608.26063657407
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
